$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of A1084 (date-formatted style s="2") down to the new A-column cells
$ws.Range("A1084").Copy() | Out-Null
$ws.Range("A1085:A1117").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the new rows 1085:1117 with OHLCV data
$ws.Cells.Item(1085, 1).Value = 45534.5
$ws.Cells.Item(1085, 2).Value = 0.10131
$ws.Cells.Item(1085, 3).Value = 0.10249
$ws.Cells.Item(1085, 4).Value = 0.09805999999999999
$ws.Cells.Item(1085, 5).Value = 0.09848
$ws.Cells.Item(1085, 6).Value = 153686611
$ws.Cells.Item(1086, 1).Value = 45534.66666666666
$ws.Cells.Item(1086, 2).Value = 0.09848999999999999
$ws.Cells.Item(1086, 3).Value = 0.10167
$ws.Cells.Item(1086, 4).Value = 0.09697
$ws.Cells.Item(1086, 5).Value = 0.10016
$ws.Cells.Item(1086, 6).Value = 147208004
$ws.Cells.Item(1087, 1).Value = 45534.83333333334
$ws.Cells.Item(1087, 2).Value = 0.10016
$ws.Cells.Item(1087, 3).Value = 0.10315
$ws.Cells.Item(1087, 4).Value = 0.10012
$ws.Cells.Item(1087, 5).Value = 0.10177
$ws.Cells.Item(1087, 6).Value = 107645920
$ws.Cells.Item(1088, 1).Value = 45535
$ws.Cells.Item(1088, 2).Value = 0.10177
$ws.Cells.Item(1088, 3).Value = 0.10239
$ws.Cells.Item(1088, 4).Value = 0.10087
$ws.Cells.Item(1088, 5).Value = 0.10126
$ws.Cells.Item(1088, 6).Value = 95799326
$ws.Cells.Item(1089, 1).Value = 45535.16666666666
$ws.Cells.Item(1089, 2).Value = 0.10127
$ws.Cells.Item(1089, 3).Value = 0.10192
$ws.Cells.Item(1089, 4).Value = 0.10097
$ws.Cells.Item(1089, 5).Value = 0.10111
$ws.Cells.Item(1089, 6).Value = 27910462
$ws.Cells.Item(1090, 1).Value = 45535.33333333334
$ws.Cells.Item(1090, 2).Value = 0.10111
$ws.Cells.Item(1090, 3).Value = 0.10152
$ws.Cells.Item(1090, 4).Value = 0.10035
$ws.Cells.Item(1090, 5).Value = 0.10125
$ws.Cells.Item(1090, 6).Value = 33504960
$ws.Cells.Item(1091, 1).Value = 45535.5
$ws.Cells.Item(1091, 2).Value = 0.10125
$ws.Cells.Item(1091, 3).Value = 0.10166
$ws.Cells.Item(1091, 4).Value = 0.10076
$ws.Cells.Item(1091, 5).Value = 0.10163
$ws.Cells.Item(1091, 6).Value = 33413862
$ws.Cells.Item(1092, 1).Value = 45535.66666666666
$ws.Cells.Item(1092, 2).Value = 0.10164
$ws.Cells.Item(1092, 3).Value = 0.10231
$ws.Cells.Item(1092, 4).Value = 0.10055
$ws.Cells.Item(1092, 5).Value = 0.10106
$ws.Cells.Item(1092, 6).Value = 49284380
$ws.Cells.Item(1093, 1).Value = 45535.83333333334
$ws.Cells.Item(1093, 2).Value = 0.10106
$ws.Cells.Item(1093, 3).Value = 0.10145
$ws.Cells.Item(1093, 4).Value = 0.10091
$ws.Cells.Item(1093, 5).Value = 0.10134
$ws.Cells.Item(1093, 6).Value = 14408827
$ws.Cells.Item(1094, 1).Value = 45536
$ws.Cells.Item(1094, 2).Value = 0.10128
$ws.Cells.Item(1094, 3).Value = 0.10153
$ws.Cells.Item(1094, 4).Value = 0.09974
$ws.Cells.Item(1094, 5).Value = 0.10028
$ws.Cells.Item(1094, 6).Value = 30930662
$ws.Cells.Item(1095, 1).Value = 45536.16666666666
$ws.Cells.Item(1095, 2).Value = 0.10027
$ws.Cells.Item(1095, 3).Value = 0.10028
$ws.Cells.Item(1095, 4).Value = 0.0985
$ws.Cells.Item(1095, 5).Value = 0.09975000000000001
$ws.Cells.Item(1095, 6).Value = 57564681
$ws.Cells.Item(1096, 1).Value = 45536.33333333334
$ws.Cells.Item(1096, 2).Value = 0.09974
$ws.Cells.Item(1096, 3).Value = 0.09984999999999999
$ws.Cells.Item(1096, 4).Value = 0.09863
$ws.Cells.Item(1096, 5).Value = 0.09936
$ws.Cells.Item(1096, 6).Value = 32497214
$ws.Cells.Item(1097, 1).Value = 45536.5
$ws.Cells.Item(1097, 2).Value = 0.09936
$ws.Cells.Item(1097, 3).Value = 0.09944
$ws.Cells.Item(1097, 4).Value = 0.09726
$ws.Cells.Item(1097, 5).Value = 0.09855999999999999
$ws.Cells.Item(1097, 6).Value = 113374135
$ws.Cells.Item(1098, 1).Value = 45536.66666666666
$ws.Cells.Item(1098, 2).Value = 0.09855999999999999
$ws.Cells.Item(1098, 3).Value = 0.09887
$ws.Cells.Item(1098, 4).Value = 0.0975
$ws.Cells.Item(1098, 5).Value = 0.0979
$ws.Cells.Item(1098, 6).Value = 23488989
$ws.Cells.Item(1099, 1).Value = 45536.83333333334
$ws.Cells.Item(1099, 2).Value = 0.09912
$ws.Cells.Item(1099, 3).Value = 0.09952
$ws.Cells.Item(1099, 4).Value = 0.09876
$ws.Cells.Item(1099, 5).Value = 0.09879
$ws.Cells.Item(1099, 6).Value = 6176878
$ws.Cells.Item(1100, 1).Value = 45537
$ws.Cells.Item(1100, 2).Value = 0.09508
$ws.Cells.Item(1100, 3).Value = 0.09606000000000001
$ws.Cells.Item(1100, 4).Value = 0.09468
$ws.Cells.Item(1100, 5).Value = 0.09522
$ws.Cells.Item(1100, 6).Value = 96705728
$ws.Cells.Item(1101, 1).Value = 45537.16666666666
$ws.Cells.Item(1101, 2).Value = 0.09521
$ws.Cells.Item(1101, 3).Value = 0.09571
$ws.Cells.Item(1101, 4).Value = 0.09408999999999999
$ws.Cells.Item(1101, 5).Value = 0.09453
$ws.Cells.Item(1101, 6).Value = 68536512
$ws.Cells.Item(1102, 1).Value = 45537.33333333334
$ws.Cells.Item(1102, 2).Value = 0.09452000000000001
$ws.Cells.Item(1102, 3).Value = 0.09802
$ws.Cells.Item(1102, 4).Value = 0.09440999999999999
$ws.Cells.Item(1102, 5).Value = 0.0973
$ws.Cells.Item(1102, 6).Value = 90853394
$ws.Cells.Item(1103, 1).Value = 45537.5
$ws.Cells.Item(1103, 2).Value = 0.0973
$ws.Cells.Item(1103, 3).Value = 0.09797
$ws.Cells.Item(1103, 4).Value = 0.0961
$ws.Cells.Item(1103, 5).Value = 0.09735000000000001
$ws.Cells.Item(1103, 6).Value = 67286158
$ws.Cells.Item(1104, 1).Value = 45537.66666666666
$ws.Cells.Item(1104, 2).Value = 0.09735000000000001
$ws.Cells.Item(1104, 3).Value = 0.09809
$ws.Cells.Item(1104, 4).Value = 0.09691
$ws.Cells.Item(1104, 5).Value = 0.09783
$ws.Cells.Item(1104, 6).Value = 42463804
$ws.Cells.Item(1105, 1).Value = 45537.83333333334
$ws.Cells.Item(1105, 2).Value = 0.09784
$ws.Cells.Item(1105, 3).Value = 0.09822
$ws.Cells.Item(1105, 4).Value = 0.09761
$ws.Cells.Item(1105, 5).Value = 0.09816999999999999
$ws.Cells.Item(1105, 6).Value = 9136695
$ws.Cells.Item(1106, 1).Value = 45538
$ws.Cells.Item(1106, 2).Value = 0.09912
$ws.Cells.Item(1106, 3).Value = 0.10081
$ws.Cells.Item(1106, 4).Value = 0.09881
$ws.Cells.Item(1106, 5).Value = 0.09957000000000001
$ws.Cells.Item(1106, 6).Value = 77084477
$ws.Cells.Item(1107, 1).Value = 45538.16666666666
$ws.Cells.Item(1107, 2).Value = 0.09958
$ws.Cells.Item(1107, 3).Value = 0.09972
$ws.Cells.Item(1107, 4).Value = 0.09901
$ws.Cells.Item(1107, 5).Value = 0.09959
$ws.Cells.Item(1107, 6).Value = 40230806
$ws.Cells.Item(1108, 1).Value = 45538.33333333334
$ws.Cells.Item(1108, 2).Value = 0.09959999999999999
$ws.Cells.Item(1108, 3).Value = 0.0997
$ws.Cells.Item(1108, 4).Value = 0.09859999999999999
$ws.Cells.Item(1108, 5).Value = 0.09962
$ws.Cells.Item(1108, 6).Value = 27094673
$ws.Cells.Item(1109, 1).Value = 45538.5
$ws.Cells.Item(1109, 2).Value = 0.09962
$ws.Cells.Item(1109, 3).Value = 0.1005
$ws.Cells.Item(1109, 4).Value = 0.09673
$ws.Cells.Item(1109, 5).Value = 0.09692000000000001
$ws.Cells.Item(1109, 6).Value = 115797619
$ws.Cells.Item(1110, 1).Value = 45538.66666666666
$ws.Cells.Item(1110, 2).Value = 0.0969
$ws.Cells.Item(1110, 3).Value = 0.09863
$ws.Cells.Item(1110, 4).Value = 0.09667000000000001
$ws.Cells.Item(1110, 5).Value = 0.09773
$ws.Cells.Item(1110, 6).Value = 47972720
$ws.Cells.Item(1111, 1).Value = 45538.83333333334
$ws.Cells.Item(1111, 2).Value = 0.09773
$ws.Cells.Item(1111, 3).Value = 0.09814000000000001
$ws.Cells.Item(1111, 4).Value = 0.09639
$ws.Cells.Item(1111, 5).Value = 0.09676999999999999
$ws.Cells.Item(1111, 6).Value = 45738801
$ws.Cells.Item(1112, 1).Value = 45539
$ws.Cells.Item(1112, 2).Value = 0.09678
$ws.Cells.Item(1112, 3).Value = 0.0975
$ws.Cells.Item(1112, 4).Value = 0.09184
$ws.Cells.Item(1112, 5).Value = 0.09581000000000001
$ws.Cells.Item(1112, 6).Value = 244041508
$ws.Cells.Item(1113, 1).Value = 45539.16666666666
$ws.Cells.Item(1113, 2).Value = 0.09581000000000001
$ws.Cells.Item(1113, 3).Value = 0.09676
$ws.Cells.Item(1113, 4).Value = 0.09449
$ws.Cells.Item(1113, 5).Value = 0.09636
$ws.Cells.Item(1113, 6).Value = 95979560
$ws.Cells.Item(1114, 1).Value = 45539.33333333334
$ws.Cells.Item(1114, 2).Value = 0.09637
$ws.Cells.Item(1114, 3).Value = 0.09697
$ws.Cells.Item(1114, 4).Value = 0.09506000000000001
$ws.Cells.Item(1114, 5).Value = 0.09556000000000001
$ws.Cells.Item(1114, 6).Value = 62215466
$ws.Cells.Item(1115, 1).Value = 45539.5
$ws.Cells.Item(1115, 2).Value = 0.09556000000000001
$ws.Cells.Item(1115, 3).Value = 0.0988
$ws.Cells.Item(1115, 4).Value = 0.09487
$ws.Cells.Item(1115, 5).Value = 0.09854
$ws.Cells.Item(1115, 6).Value = 136433380
$ws.Cells.Item(1116, 1).Value = 45539.66666666666
$ws.Cells.Item(1116, 2).Value = 0.09853000000000001
$ws.Cells.Item(1116, 3).Value = 0.09951
$ws.Cells.Item(1116, 4).Value = 0.09759
$ws.Cells.Item(1116, 5).Value = 0.09798999999999999
$ws.Cells.Item(1116, 6).Value = 91680872
$ws.Cells.Item(1117, 1).Value = 45539.83333333334
$ws.Cells.Item(1117, 2).Value = 0.09805
$ws.Cells.Item(1117, 3).Value = 0.09812
$ws.Cells.Item(1117, 4).Value = 0.0979
$ws.Cells.Item(1117, 5).Value = 0.09795
$ws.Cells.Item(1117, 6).Value = 157266
